$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title: Facture -> Invoice
$ws.Range("A1").Value = "Invoice"

# Header row translations (row 4)
$ws.Range("B4").Value = "Weight/m (kg)"
$ws.Range("C4").Value = "Stock Length"
$ws.Range("D4").Value = "QTY"
$ws.Range("E4").Value = "Total Weight (kg)"
$ws.Range("F4").Value = "Percentage"
$ws.Range("G4").Value = "Total Price"

# "Prix Total" / "Total Price" column (G5:G38) holds numeric-looking text
# (e.g. "322.000"). Force the cells to Text format first so Excel keeps
# them as plain text (shared strings) instead of auto-converting to numbers.
$priceCol = $ws.Range("G5:G38")
$priceCol.NumberFormat = "@"

$ws.Range("G5").Value = "322.000"
$ws.Range("G6").Value = "48.000"
$ws.Range("G7").Value = "192.000"
$ws.Range("G8").Value = "224.000"
$ws.Range("G9").Value = "32.000"
$ws.Range("G10").Value = "576.000"
$ws.Range("G11").Value = "48.000"
$ws.Range("G12").Value = "240.000"
$ws.Range("G13").Value = "26.000"
$ws.Range("G14").Value = "124.000"
$ws.Range("G15").Value = "13908.000"
$ws.Range("G16").Value = "192.000"
$ws.Range("G17").Value = "24.000"
$ws.Range("G18").Value = "28.000"
$ws.Range("G19").Value = "42.000"
$ws.Range("G20").Value = "46.000"
$ws.Range("G21").Value = "25.000"
$ws.Range("G22").Value = "25.000"
$ws.Range("G23").Value = "180.000"
$ws.Range("G24").Value = "36.000"
$ws.Range("G25").Value = "38.000"
$ws.Range("G26").Value = "26.000"
$ws.Range("G27").Value = "1176.000"
$ws.Range("G28").Value = "408.000"
$ws.Range("G29").Value = "8997.000"
$ws.Range("G30").Value = "9233.000"
$ws.Range("G31").Value = "982.000"
$ws.Range("G32").Value = "3179.000"
$ws.Range("G33").Value = "2148.000"
$ws.Range("G34").Value = "756.000"
$ws.Range("G35").Value = "273.000"
$ws.Range("G36").Value = "341.000"
$ws.Range("G37").Value = "43895.000"

# Row 38: adjusted total label translation and value
$ws.Range("A38").Value = "ADJUSTED TOTAL (+{0}%) (12.0%)"
$ws.Range("G38").Value = "49162.400"
